$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" to H1, matching the formatting of the other header cells
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add data values for the new "Save" column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
